# MXB103StatementOfContributionTemplate.xlsx — "Add files via upload"
#
# The author re-uploaded the workbook after updating the contribution
# summary for the "Jay Choi" row: Sheet1!C5 grows from
#   "Analysis: 3,4"
# to
#   "Analysis: 3,4, Report"
# (everything else in the published diff — locale-flavoured font names,
# theme/style display names, the fileVersion/absPath build stamps, the
# phonetic-guide hints, and the font-metric-driven row height / column
# width drift — is a side effect of that save happening in a different,
# Korean-locale copy of Excel, not a deliberate content edit.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The actual edit: update the contribution note for Jay Choi.
$ws.Range("C5").Value = "Analysis: 3,4, Report"

# The re-save also left the cursor parked on C15.
[void]$ws.Range("C15").Select()
